$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table held yearly data for 2005-2020 in rows 2-17.
# It now needs to hold 2010-2021 in rows 2-13:
#  - drop the five oldest years (2005-2009), shifting 2010-2020 up to rows 2-12
#  - append a new row for 2021 at row 13

# Remove rows 2-6 (2005年..2009年); this shifts 2010年..2020年 up into rows 2-12
$ws.Range("A2:A6").EntireRow.Delete()

# Add the new 2021 row at row 13, copying the formatting from row 12 (the
# previous last row) so the year label keeps the same style.
$ws.Range("A12").Copy($ws.Range("A13"))

$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 155.6589
$ws.Range("C13").Value = 2600.2281
$ws.Range("D13").Value = 3637.6865
$ws.Range("E13").Value = 1125.3502
$ws.Range("F13").Value = 7.8029
$ws.Range("G13").Value = 7168.3363
$ws.Range("H13").Value = 1689.7343
$ws.Range("I13").Value = 12867.2589
$ws.Range("J13").Value = 460.7987
